$d = $word.ActiveDocument

$replacements = @(
    @("67×35=", "51×82="),
    @("25×85=", "78×37="),
    @("35×52=", "80×48="),
    @("50×49=", "38×29="),
    @("83×27=", "58×78="),
    @("53×21=", "75×35="),
    @("73×52=", "77×73="),
    @("87×86=", "11×47="),
    @("26×42=", "94×71="),
    @("71×91=", "95×29="),
    @("34×99=", "46×89="),
    @("50×23=", "12×30="),
    @("37×28=", "21×71="),
    @("28×94=", "77×79="),
    @("64×77=", "86×17="),
    @("29×98=", "19×51="),
    @("66×93=", "57×73="),
    @("88×82=", "93×60="),
    @("76×37=", "86×58="),
    @("46×99=", "62×84="),
    @("85×19=", "48×42="),
    @("28×81=", "17×41="),
    @("40×46=", "59×95="),
    @("95×45=", "47×86="),
    @("20×48=", "65×46=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
